$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.3754574605433065"
$ws.Range("E2").Value = [double]"0.3754574605433065"

# Row 3
$ws.Range("D3").Value = [double]"0.3172054871708445"
$ws.Range("E3").Value = [double]"0.3172054871708445"

# Row 4
$ws.Range("D4").Value = [double]"0.002998188921254697"
$ws.Range("E4").Value = [double]"0.002998188921254697"

# Row 5
$ws.Range("D5").Value = [double]"7.48909733850698E-15"
$ws.Range("E5").Value = [double]"7.48909733850698E-15"

# Row 6
$ws.Range("D6").Value = [double]"0.6817466841128406"
$ws.Range("E6").Value = [double]"0.6817466841128406"

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = [double]"0.1004979671864193"
$ws.Range("E7").Value = [double]"0.8995020328135808"

# Row 8
$ws.Range("D8").Value = [double]"0.9999997473473775"
$ws.Range("E8").Value = [double]"2.526526224899683E-07"

# Row 9
$ws.Range("D9").Value = [double]"0.5775780280966436"
$ws.Range("E9").Value = [double]"0.4224219719033564"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"0.00761718701365999"
$ws.Range("E10").Value = [double]"0.99238281298634"

# Row 11
$ws.Range("D11").Value = [double]"0.9956277316222051"
$ws.Range("E11").Value = [double]"0.004372268377794919"
$ws.Range("F11").Value = [double]"0.9728466272354126"
$ws.Range("G11").Value = [double]"0.7"

# Row 12
$ws.Range("D12").Value = [double]"0.0001406487130089631"
$ws.Range("E12").Value = [double]"0.0001406487130089631"

# Row 13
$ws.Range("D13").Value = [double]"0.09909819165592877"
$ws.Range("E13").Value = [double]"0.09909819165592877"

# Row 14
$ws.Range("D14").Value = [double]"0.149456476801289"
$ws.Range("E14").Value = [double]"0.149456476801289"

# Row 15
$ws.Range("D15").Value = [double]"7.112495171131203E-14"
$ws.Range("E15").Value = [double]"7.112495171131203E-14"

# Row 16
$ws.Range("D16").Value = [double]"0.5834470176022701"
$ws.Range("E16").Value = [double]"0.5834470176022701"

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = [double]"0.0002866432735851267"
$ws.Range("E17").Value = [double]"0.9997133567264149"

# Row 18
$ws.Range("D18").Value = [double]"0.99999999996588"
$ws.Range("E18").Value = [double]"3.412004012659509E-11"

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"0.02339736043789168"
$ws.Range("E19").Value = [double]"0.9766026395621084"

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = [double]"0.9961457873023433"
$ws.Range("E21").Value = [double]"0.003854212697656734"
$ws.Range("F21").Value = [double]"1.305838823318481"
$ws.Range("G21").Value = [double]"0.7"

$wb.Save()
